$wb = $excel.ActiveWorkbook

# Target values per sheet: row -> @(C, D, E, F, G, I)
$sheetData = @{}

$sheetData["Tarantula"] = @{
    2 = @(-1, 3, 0, 3, 3, 3)
    3 = @(-1, 7, 0, 7, 7, 16)
    4 = @(-1, 3, 0, 3, 3, 3)
    5 = @(-1, 4, 0, 4, 4, 6)
    6 = @(-1, 6, 0, 6, 6, 10)
    7 = @(-1, 3, 0, 3, 3, 28)
    8 = @(-1, 7, 0, 7, 7, 13)
    9 = @(-1, 12, 0, 12, 12, 13)
    10 = @(-1, 3, 0, 3, 3, 11)
    11 = @(-1, 13, 0, 13, 13, 3)
    12 = @(-1, 1, 0, 1, 1, 5)
    13 = @(-1, 4, 0, 4, 4, 6)
    14 = @(-1, 3, 0, 3, 3, 3)
    15 = @(-1, 19, 0, 19, 19, 4)
    16 = @(-1, 3, 0, 3, 3, 2)
    17 = @(-1, 5, 0, 5, 5, 6)
    18 = @(-1, 1, 0, 1, 1, 1)
    19 = @(-1, 25, 0, 25, 25, 3)
    20 = @(-1, 5, 0, 5, 5, 3)
    21 = @(-1, 5, 0, 5, 5, 11)
    22 = @(-1, 24, 0, 24, 24, 2)
    23 = @(-1, 3, 0, 3, 3, 4)
    24 = @(-1, 4, 0, 4, 4, 4)
    25 = @(-1, 3, 0, 3, 3, 13)
    26 = @(-1, 20, 0, 20, 20, 3)
    27 = @(-1, 3, 0, 3, 3, 13)
    28 = @(-1, 5, 0, 5, 5, 10)
    29 = @(-1, 3, 0, 3, 3, 12)
    30 = @(-1, 20, 0, 20, 20, 3)
    31 = @(-1, 2, 0, 2, 2, 3)
    32 = @(-1, 6, 0, 6, 6, 7)
    33 = @(-1, 5, 0, 5, 5, 5)
    34 = @(-1, 5, 0, 5, 5, 13)
}

$sheetData["Ochiai"] = @{
    2 = @(-1, 3, 0, 3, 3, 3)
    3 = @(-1, 3, 0, 3, 3, 2)
    4 = @(-1, 3, 0, 3, 3, 3)
    5 = @(-1, 4, 0, 4, 4, 4)
    6 = @(-1, 3, 0, 3, 3, 3)
    7 = @(-1, 3, 0, 3, 3, 3)
    8 = @(-1, 3, 0, 3, 3, 3)
    9 = @(-1, 3, 0, 3, 3, 3)
    10 = @(-1, 3, 0, 3, 3, 2)
    11 = @(-1, 13, 0, 13, 13, 3)
    12 = @(-1, 1, 0, 1, 1, 5)
    13 = @(-1, 4, 0, 4, 4, 4)
    14 = @(-1, 3, 0, 3, 3, 3)
    15 = @(-1, 19, 0, 19, 19, 4)
    16 = @(-1, 3, 0, 3, 3, 2)
    17 = @(-1, 5, 0, 5, 5, 3)
    18 = @(-1, 1, 0, 1, 1, 1)
    19 = @(-1, 25, 0, 25, 25, 2)
    20 = @(-1, 5, 0, 5, 5, 3)
    21 = @(-1, 5, 0, 5, 5, 2)
    22 = @(-1, 31, 0, 31, 31, 14)
    23 = @(-1, 3, 0, 3, 3, 3)
    24 = @(-1, 4, 0, 4, 4, 4)
    25 = @(-1, 3, 0, 3, 3, 3)
    26 = @(-1, 20, 0, 20, 20, 2)
    27 = @(-1, 3, 0, 3, 3, 3)
    28 = @(-1, 5, 0, 5, 5, 9)
    29 = @(-1, 3, 0, 3, 3, 2)
    30 = @(-1, 20, 0, 20, 20, 2)
    31 = @(-1, 2, 0, 2, 2, 2)
    32 = @(-1, 6, 0, 6, 6, 4)
    33 = @(-1, 5, 0, 5, 5, 5)
    34 = @(-1, 5, 0, 5, 5, 3)
}

$sheetData["Op2"] = @{
    2 = @(-1, 3, 0, 3, 3, 3)
    3 = @(-1, 3, 0, 3, 3, 2)
    4 = @(-1, 3, 0, 3, 3, 3)
    5 = @(-1, 4, 0, 4, 4, 4)
    6 = @(-1, 3, 0, 3, 3, 3)
    7 = @(-1, 3, 0, 3, 3, 3)
    8 = @(-1, 3, 0, 3, 3, 3)
    9 = @(-1, 3, 0, 3, 3, 3)
    10 = @(-1, 3, 0, 3, 3, 2)
    11 = @(-1, 13, 0, 13, 13, 3)
    12 = @(-1, 1, 0, 1, 1, 5)
    13 = @(-1, 4, 0, 4, 4, 4)
    14 = @(-1, 3, 0, 3, 3, 3)
    15 = @(-1, 19, 0, 19, 19, 4)
    16 = @(-1, 3, 0, 3, 3, 2)
    17 = @(-1, 5, 0, 5, 5, 3)
    18 = @(-1, 1, 0, 1, 1, 1)
    19 = @(-1, 25, 0, 25, 25, 2)
    20 = @(-1, 5, 0, 5, 5, 3)
    21 = @(-1, 5, 0, 5, 5, 2)
    22 = @(-1, 31, 0, 31, 31, 28)
    23 = @(-1, 3, 0, 3, 3, 3)
    24 = @(-1, 4, 0, 4, 4, 4)
    25 = @(-1, 3, 0, 3, 3, 3)
    26 = @(-1, 20, 0, 20, 20, 2)
    27 = @(-1, 3, 0, 3, 3, 3)
    28 = @(-1, 5, 0, 5, 5, 9)
    29 = @(-1, 3, 0, 3, 3, 2)
    30 = @(-1, 20, 0, 20, 20, 2)
    31 = @(-1, 2, 0, 2, 2, 2)
    32 = @(-1, 6, 0, 6, 6, 4)
    33 = @(-1, 5, 0, 5, 5, 5)
    34 = @(-1, 5, 0, 5, 5, 3)
}

$sheetData["Barinel"] = @{
    2 = @(-1, 3, 0, 3, 3, 3)
    3 = @(-1, 7, 0, 7, 7, 16)
    4 = @(-1, 3, 0, 3, 3, 3)
    5 = @(-1, 4, 0, 4, 4, 6)
    6 = @(-1, 6, 0, 6, 6, 10)
    7 = @(-1, 3, 0, 3, 3, 28)
    8 = @(-1, 7, 0, 7, 7, 13)
    9 = @(-1, 12, 0, 12, 12, 13)
    10 = @(-1, 3, 0, 3, 3, 11)
    11 = @(-1, 13, 0, 13, 13, 3)
    12 = @(-1, 1, 0, 1, 1, 5)
    13 = @(-1, 4, 0, 4, 4, 6)
    14 = @(-1, 3, 0, 3, 3, 3)
    15 = @(-1, 19, 0, 19, 19, 4)
    16 = @(-1, 3, 0, 3, 3, 2)
    17 = @(-1, 5, 0, 5, 5, 6)
    18 = @(-1, 1, 0, 1, 1, 1)
    19 = @(-1, 25, 0, 25, 25, 3)
    20 = @(-1, 5, 0, 5, 5, 3)
    21 = @(-1, 5, 0, 5, 5, 11)
    22 = @(-1, 24, 0, 24, 24, 2)
    23 = @(-1, 3, 0, 3, 3, 4)
    24 = @(-1, 4, 0, 4, 4, 4)
    25 = @(-1, 3, 0, 3, 3, 13)
    26 = @(-1, 20, 0, 20, 20, 3)
    27 = @(-1, 3, 0, 3, 3, 13)
    28 = @(-1, 5, 0, 5, 5, 10)
    29 = @(-1, 3, 0, 3, 3, 12)
    30 = @(-1, 20, 0, 20, 20, 3)
    31 = @(-1, 2, 0, 2, 2, 3)
    32 = @(-1, 6, 0, 6, 6, 7)
    33 = @(-1, 5, 0, 5, 5, 5)
    34 = @(-1, 5, 0, 5, 5, 13)
}

$sheetData["Dstar"] = @{
    2 = @(-1, 3, 0, 3, 3, 3)
    3 = @(-1, 3, 0, 3, 3, 2)
    4 = @(-1, 95, 0, 95, 95, 3)
    5 = @(-1, 4, 0, 4, 4, 4)
    6 = @(-1, 3, 0, 3, 3, 3)
    7 = @(-1, 3, 0, 3, 3, 3)
    8 = @(-1, 3, 0, 3, 3, 3)
    9 = @(-1, 3, 0, 3, 3, 3)
    10 = @(-1, 3, 0, 3, 3, 2)
    11 = @(-1, 13, 0, 13, 13, 3)
    12 = @(-1, 1, 0, 1, 1, 5)
    13 = @(-1, 4, 0, 4, 4, 4)
    14 = @(-1, 3, 0, 3, 3, 3)
    15 = @(-1, 19, 0, 19, 19, 4)
    16 = @(-1, 15, 0, 15, 15, 2)
    17 = @(-1, 5, 0, 5, 5, 3)
    18 = @(-1, 1, 0, 1, 1, 1)
    19 = @(-1, 25, 0, 25, 25, 2)
    20 = @(-1, 5, 0, 5, 5, 3)
    21 = @(-1, 5, 0, 5, 5, 2)
    22 = @(-1, 31, 0, 31, 31, 14)
    23 = @(-1, 3, 0, 3, 3, 3)
    24 = @(-1, 4, 0, 4, 4, 4)
    25 = @(-1, 3, 0, 3, 3, 3)
    26 = @(-1, 20, 0, 20, 20, 2)
    27 = @(-1, 3, 0, 3, 3, 3)
    28 = @(-1, 5, 0, 5, 5, 9)
    29 = @(-1, 3, 0, 3, 3, 2)
    30 = @(-1, 20, 0, 20, 20, 2)
    31 = @(-1, 2, 0, 2, 2, 2)
    32 = @(-1, 6, 0, 6, 6, 4)
    33 = @(-1, 5, 0, 5, 5, 5)
    34 = @(-1, 5, 0, 5, 5, 13)
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetData[$sheetName]
    foreach ($row in $rows.Keys) {
        $vals = $rows[$row]
        $ws.Cells.Item($row, 3).Value = $vals[0]   # C
        $ws.Cells.Item($row, 4).Value = $vals[1]   # D
        $ws.Cells.Item($row, 5).Value = $vals[2]   # E
        $ws.Cells.Item($row, 6).Value = $vals[3]   # F
        $ws.Cells.Item($row, 7).Value = $vals[4]   # G
        $ws.Cells.Item($row, 9).Value = $vals[5]   # I
    }
}